$d = $word.ActiveDocument

# --- Insertion 1: add "signature.page.scss" line after
#     "src/app/members/signature/signature.page.html" and before
#     "src/app/members/signature/signature.page.ts" ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("signature.page.html", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter([char]11 + "src/app/members/signature/signature.page.scss")
}

# --- Insertion 2: add "login.page.scss" line after
#     "src/app/public/login/login.page.html" and before
#     "src/app/public/login/login.page.ts" ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("login.page.html", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter([char]11 + "src/app/public/login/login.page.scss")
}
